# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1375
    $ws.Range("F3").Value = 2146
    $ws.Range("F4").Value = 303
    $ws.Range("F6").Value = 6398
    $ws.Range("F7").Value = 274
}
